$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.054.90"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.566.51"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.46"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.13"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.566.97"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.77"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "27.059.82"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.91"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.65"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.19"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "153.79"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.23"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").Value = "1.423.69"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +12.76%  "
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.814"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.61"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "1.704.07"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.84"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  +0.19%  "
